$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update OLBG_Confidence values (column D)
$ws.Range("D2").Value = 96
$ws.Range("D3").Value = 95
$ws.Range("D4").Value = 98

# Fill in Odds values (column F) for rows previously blank
$ws.Range("F6").Value = 1.25
$ws.Range("F7").Value = 1.25
$ws.Range("F8").Value = 1.33
$ws.Range("F9").Value = 2.3
